$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style (bold, centered, bordered) from H1 into the two new header cells
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$iValues = @(6,7,6,7,8,8,9,3,9,9,6,6,6,6,8,9,7,7,8,7,8,8,5,8,5,1,7,8,7,1,7,2,4,8,8,8,8,5,4,8,8,7,6,3,8,6,6,7,7,6,7,6,7,5,7,8,6,6,2,2,5,4)
$jValues = @(7,8,7,8,8,8,9,3,9,9,7,6,6,6,8,9,8,8,8,8,8,8,6,8,6,3,7,8,7,2,7,5,5,8,8,8,8,6,5,8,8,7,6,4,8,7,6,7,8,7,8,6,7,6,7,9,7,7,3,3,5,4)

for ($n = 0; $n -lt $iValues.Length; $n++) {
    $row = $n + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$n]
    $ws.Cells.Item($row, 10).Value = $jValues[$n]
}

